$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.303.06'
$ws.Range("E2").Value = '  +9.29%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.229.89'
$ws.Range("E3").Value = '  +4.05%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '399.36'
$ws.Range("E5").Value = '  +3.51%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '110.69'
$ws.Range("E6").Value = '  +6.36%  '

# Row 7
$ws.Range("E7").Value = '  +2.75%  '

# Row 8
$ws.Range("E8").Value = '  -0.03%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.625'
$ws.Range("E9").Value = '  +6.66%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.46'
$ws.Range("E10").Value = '  +6.12%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0908'
$ws.Range("E11").Value = '  +5.87%  '

# Row 12
$ws.Range("E12").Value = '  +2.20%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.739.28'
$ws.Range("E13").Value = '  +3.99%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.13'
$ws.Range("E14").Value = '  +3.87%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.14'
$ws.Range("E15").Value = '  +3.17%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.238.90'
$ws.Range("E16").Value = '  +4.32%  '

# Row 17
$ws.Range("E17").Value = '  +5.66%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.70'
$ws.Range("E18").Value = '  -2.01%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '56.127.66'
$ws.Range("E19").Value = '  +8.81%  '

# Row 20
$ws.Range("E20").Value = '  +2.01%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000102'
$ws.Range("E21").Value = '  +6.37%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.08'
$ws.Range("E22").Value = '  +4.45%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '304.31'
$ws.Range("E23").Value = '  +14.04%  '

# Row 25
$ws.Range("E25").Value = '  +2.54%  '

# Row 26
$ws.Range("E26").Value = '  +1.47%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '28.25'
$ws.Range("E27").Value = '  +4.38%  '

# Row 28
$ws.Range("E28").Value = '  +3.32%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.173'
$ws.Range("E29").Value = '  +4.16%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.09%  '

# Row 31
$ws.Range("E31").Value = '  +4.58%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.21'
$ws.Range("E32").Value = '  +7.48%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0494'
$ws.Range("E33").Value = '  +3.00%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '36.32'
$ws.Range("E34").Value = '  +2.26%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.12'
$ws.Range("E35").Value = '  +2.45%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '51.39'
$ws.Range("E36").Value = '  +2.70%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.11'
$ws.Range("E37").Value = '  +23.94%  '

# Row 38
$ws.Range("E38").Value = '  +0.03%  '

# Row 39
$ws.Range("E39").Value = '  +4.06%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '134.94'
$ws.Range("E40").Value = '  +4.61%  '

# Row 41
$ws.Range("E41").Value = '  +3.12%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.04'
$ws.Range("E42").Value = '  +6.66%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.25'
$ws.Range("E43").Value = '  +3.89%  '

# Row 44
$ws.Range("E44").Value = '  +3.45%  '

# Row 45
$ws.Range("E45").Value = '  -2.70%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.32'
$ws.Range("E46").Value = '  +0.73%  '

# Row 47
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("E47").Value = '  +2.48%  '

# Row 48
$ws.Range("B48").Value = 'ThetaToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.13'
$ws.Range("E48").Value = '  +47.55%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.48'
$ws.Range("E49").Value = '  -1.79%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.138.51'
$ws.Range("E50").Value = '  +3.02%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0362'
$ws.Range("E51").Value = '  +9.86%  '
